$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 98
$ws.Range("F3").Value = 124
$ws.Range("F4").Value = 635
$ws.Range("F5").Value = 364
$ws.Range("F6").Value = 549
$ws.Range("F7").Value = 1539
$ws.Range("F9").Value = 11723
$ws.Range("F13").Value = 2102
$ws.Range("F15").Value = 37
$ws.Range("F19").Value = 1195
$ws.Range("F20").Value = 171
$ws.Range("F22").Value = 737
$ws.Range("F23").Value = 658
$ws.Range("F24").Value = 271
$ws.Range("F26").Value = 728
$ws.Range("F27").Value = 3597
$ws.Range("F28").Value = 3597
$ws.Range("F30").Value = 812
$ws.Range("F34").Value = 989
$ws.Range("F35").Value = 36
$ws.Range("F37").Value = 255
$ws.Range("F41").Value = 4424
$ws.Range("F42").Value = 5471
$ws.Range("F46").Value = 269
$ws.Range("F47").Value = 66
$ws.Range("F48").Value = 29
$ws.Range("F50").Value = 102

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 4155
$ws.Range("F5").Value = 90
$ws.Range("F11").Value = 685

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 752
$ws.Range("F3").Value = 422
$ws.Range("F4").Value = 63

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 752
$ws.Range("F3").Value = 422
$ws.Range("F4").Value = 63
$ws.Range("F5").Value = 98
$ws.Range("F6").Value = 124
$ws.Range("F7").Value = 635
$ws.Range("F8").Value = 364
$ws.Range("F9").Value = 549
$ws.Range("F10").Value = 1539
$ws.Range("F11").Value = 11723
$ws.Range("F15").Value = 2102
$ws.Range("F17").Value = 37
$ws.Range("F20").Value = 1195
$ws.Range("F21").Value = 171
$ws.Range("F23").Value = 4155
$ws.Range("F25").Value = 271
$ws.Range("F26").Value = 728
$ws.Range("F27").Value = 3597
$ws.Range("F29").Value = 90
$ws.Range("F31").Value = 812
$ws.Range("F32").Value = 989
$ws.Range("F33").Value = 36
$ws.Range("F35").Value = 255
$ws.Range("F38").Value = 4424
$ws.Range("F42").Value = 269
$ws.Range("F45").Value = 66
$ws.Range("F46").Value = 29
$ws.Range("F50").Value = 102

